$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it used to wrap the first
#    few paragraphs of the document).
# ---------------------------------------------------------------------
$bms = $d.Bookmarks
if ($bms.Exists("_GoBack")) {
    $bms.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Insert a brand-new list paragraph right after the paragraph that
#    contains "Pleier dere å gi pasienten program ...".
#    InsertParagraphAfter() duplicates the pPr/rPr of the current
#    paragraph (numbering, spacing, fonts, ...), which already matches
#    what the diff expects for the new paragraph.
# ---------------------------------------------------------------------
$srcRange = $d.Content
$srcRange.Find.Execute("Pleier dere å gi pasienten program som de må trene på hjemme mellom hver time?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$srcRange.InsertParagraphAfter()

$srcPara = $srcRange.Paragraphs(1)
$newPara = $srcPara.Next()
$newParaRange = $newPara.Range

# Range covering everything in the new paragraph except its trailing
# paragraph mark.
$textRange = $d.Range($newParaRange.Start, $newParaRange.End - 1)

# Temporarily add a placeholder character so that the position right
# after the sentence is NOT the very edge of the paragraph; this lets
# us drop a properly collapsed bookmark there without the engine
# widening it to cover neighbouring paragraphs or splitting the run.
$textRange.Text = "Hvordan legger dere opp treningen? Og hvordan legger dere opp trening hos dere vs. trening utenfor behandlingen?X"

# ---------------------------------------------------------------------
# 3. Re-locate the inserted sentence (now that the document has
#    changed) and drop the "_GoBack" bookmark around the placeholder
#    character sitting right at its end.
# ---------------------------------------------------------------------
$bmFindRange = $d.Content
$bmFindRange.Find.Execute("Hvordan legger dere opp treningen? Og hvordan legger dere opp trening hos dere vs. trening utenfor behandlingen?X", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$placeholderRange = $d.Range($bmFindRange.End - 1, $bmFindRange.End)
$d.Bookmarks.Add("_GoBack", $placeholderRange)

# Remove the placeholder character again; the bookmark collapses to a
# zero-length range right after the sentence, exactly where it belongs.
$placeholderRange2 = $d.Range($bmFindRange.End - 1, $bmFindRange.End)
$placeholderRange2.Text = ""
